$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("#system")

# 1. Insert a new column before column N (14) -- this shifts the cell
#    contents of N:AC to O:AD, leaving A:M untouched. (Defined names are
#    NOT auto-adjusted by this engine, so every affected named range is
#    fixed up by hand further below.)
$ws.Columns.Item(14).Insert()

# 2. Populate the new column N with the "localdb" command group.
$ws.Range("N1").Value = "localdb"
$ws.Range("N2").Value = "cloneTable(var,source,target)"
$ws.Range("N3").Value = "dropTables(var,tables)"
$ws.Range("N4").Value = "exportCSV(sql,output)"
$ws.Range("N5").Value = "importRecords(var,sourceDb,sql,table)"
$ws.Range("N6").Value = "purge(var)"
$ws.Range("N7").Value = "runSQLs(var,sqls)"

# 3. Insert "localdb" into the alphabetical "target" list in column A,
#    between "json" (row 13) and "macro" (row 14). A whole-row/whole-
#    column Insert() is used for step 1 above, but a single-cell
#    Range.Insert() here would incorrectly ripple into neighbouring
#    columns on this engine, so the one-column shift is done by hand:
#    walk the tail of the list upward, copying each cell into the one
#    below it, then drop "localdb" into the now-vacant row 14.
for ($r = 29; $r -ge 14; $r--) {
    $src = $ws.Range("A" + $r)
    $dst = $ws.Range("A" + ($r + 1))
    $dst.Value = $src.Text
}
$ws.Range("A14").Value = "localdb"

# 4. Fix up every defined name whose range shifted because of the
#    column insert, add the brand-new "localdb" name, and grow "target"
#    to cover the newly-inserted row.
$wb.Names.Add("localdb", '=''#system''!$N$2:$N$7')
$wb.Names.Item("macro").RefersTo      = '=''#system''!$O$2:$O$4'
$wb.Names.Item("mail").RefersTo       = '=''#system''!$P$2:$P$2'
$wb.Names.Item("number").RefersTo     = '=''#system''!$Q$2:$Q$16'
$wb.Names.Item("pdf").RefersTo        = '=''#system''!$R$2:$R$16'
$wb.Names.Item("rdbms").RefersTo      = '=''#system''!$S$2:$S$7'
$wb.Names.Item("redis").RefersTo      = '=''#system''!$T$2:$T$10'
$wb.Names.Item("sms").RefersTo        = '=''#system''!$U$2:$U$2'
$wb.Names.Item("sound").RefersTo      = '=''#system''!$V$2:$V$5'
$wb.Names.Item("ssh").RefersTo        = '=''#system''!$W$2:$W$9'
$wb.Names.Item("step").RefersTo       = '=''#system''!$X$2:$X$4'
$wb.Names.Item("target").RefersTo     = '=''#system''!$A$2:$A$30'
$wb.Names.Item("web").RefersTo        = '=''#system''!$Y$2:$Y$127'
$wb.Names.Item("webalert").RefersTo   = '=''#system''!$Z$2:$Z$8'
$wb.Names.Item("webcookie").RefersTo  = '=''#system''!$AA$2:$AA$8'
$wb.Names.Item("ws").RefersTo         = '=''#system''!$AB$2:$AB$17'
$wb.Names.Item("ws.async").RefersTo   = '=''#system''!$AC$2:$AC$8'
$wb.Names.Item("xml").RefersTo        = '=''#system''!$AD$2:$AD$21'
